# OLX Monitor update 2026-02-21 07:54
# Appends 8 new listing rows (99-106) to the "PODSUMOWANIE" sheet, which
# holds the detailed per-offer table (columns A-H) below the summary block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PODSUMOWANIE")

$rows = @(
  @{ Row=99;  A="2026-02-21 07:54:59"; B="poqui";           C="Mieszkanie z KLIMATYZACJĄ 5 minut od UMCS, UP, KUL - Długosza";                         D=2049;  E="19.12.2025"; F=63;  FStyle="s15"; G="https://www.olx.pl/d/oferta/mieszkanie-z-klimatyzacja-5-minut-od-umcs-up-kul-dlugosza-CID3-ID18KAEc.html";                           H="mieszkanie-z-klimatyzacja-5-minut-od-umcs-up-kul-dlugosza-CID3-ID18KAEc" },
  @{ Row=100; A="2026-02-21 07:54:59"; B="poqui";           C="Świeżo wykończone mieszkanie z dużym balkonem - Ponikwoda";                             D=2299;  E="19.01.2026"; F=32;  FStyle="s14"; G="https://www.olx.pl/d/oferta/swiezo-wykonczone-mieszkanie-z-duzym-balkonem-ponikwoda-CID3-ID1951OR.html";                             H="swiezo-wykonczone-mieszkanie-z-duzym-balkonem-ponikwoda-CID3-ID1951OR" },
  @{ Row=101; A="2026-02-21 07:54:59"; B="poqui";           C="Kawalerka po remoncie z funkcjonalną antresolą - ul. Jana Sawy";                        D=2499;  E="28.10.2025"; F=115; FStyle="s15"; G="https://www.olx.pl/d/oferta/kawalerka-po-remoncie-z-funkcjonalna-antresola-ul-jana-sawy-CID3-ID183ger.html";                         H="kawalerka-po-remoncie-z-funkcjonalna-antresola-ul-jana-sawy-CID3-ID183ger" },
  @{ Row=102; A="2026-02-21 07:54:59"; B="poqui";           C="Przytulny pokój blisko Politechniki – ul. Przytulna";                                    D=549;   E="10.10.2025"; F=133; FStyle="s15"; G="https://www.olx.pl/d/oferta/przytulny-pokoj-blisko-politechniki-ul-przytulna-CID3-ID17NeTz.html";                                   H="przytulny-pokoj-blisko-politechniki-ul-przytulna-CID3-ID17NeTz" },
  @{ Row=103; A="2026-02-21 07:54:59"; B="pokojewlublinie"; C="WOLNY OD ZARAZ! Pokój jedynka, ul. Romanowskiego 58";                                    D=58640; E="11.08.2025"; F=193; FStyle="s15"; G="https://www.olx.pl/d/oferta/wolny-od-zaraz-pokoj-jedynka-ul-romanowskiego-58-CID3-ID16ZeYm.html";                                   H="wolny-od-zaraz-pokoj-jedynka-ul-romanowskiego-58-CID3-ID16ZeYm" },
  @{ Row=104; A="2026-02-21 07:54:59"; B="pokojewlublinie"; C="WOLNY OD ZARAZ! Super lokalizacja, blisko centrum, ul. Paganiniego 12";                  D=12640; E="19.01.2026"; F=32;  FStyle="s14"; G="https://www.olx.pl/d/oferta/wolny-od-zaraz-super-lokalizacja-blisko-centrum-ul-paganiniego-12-CID3-ID195dLc.html";                  H="wolny-od-zaraz-super-lokalizacja-blisko-centrum-ul-paganiniego-12-CID3-ID195dLc" },
  @{ Row=105; A="2026-02-21 07:54:59"; B="dawnypatron";     C="Ładny pokój jednoosobowy. Wynajmę duży pokój w centrum. ul Niecała 4.";                  D=730;   E="20.09.2024"; F=518; FStyle="s15"; G="https://www.olx.pl/d/oferta/ladny-pokoj-jednoosobowy-wynajme-duzy-pokoj-w-centrum-ul-niecala-4-CID3-ID122jPM.html";                  H="ladny-pokoj-jednoosobowy-wynajme-duzy-pokoj-w-centrum-ul-niecala-4-CID3-ID122jPM" },
  @{ Row=106; A="2026-02-21 07:54:59"; B="dawnypatron";     C="Mam do wynajęcia pokój dla os. pracującej lub studenta. Narutowicza 14";                 D=14690; E="05.12.2025"; F=77;  FStyle="s15"; G="https://www.olx.pl/d/oferta/mam-do-wynajecia-pokoj-dla-os-pracujacej-lub-studenta-narutowicza-14-CID3-ID18ySfv.html";                 H="mam-do-wynajecia-pokoj-dla-os-pracujacej-lub-studenta-narutowicza-14-CID3-ID18ySfv" }
)

# --- 1. Write the raw values --------------------------------------------
# Text columns are force-typed to Text first (NumberFormat "@") so Excel's
# locale-aware "smart" parser does not silently reinterpret an unambiguous
# looking date string (e.g. "10.10.2025", "11.08.2025", "05.12.2025") as a
# real date serial - every value in this table is plain text in the source.
foreach ($r in $rows) {
  $ws.Cells.Item($r.Row, 1).NumberFormat = "@"
  $ws.Cells.Item($r.Row, 1).Value = $r.A

  $ws.Cells.Item($r.Row, 2).NumberFormat = "@"
  $ws.Cells.Item($r.Row, 2).Value = $r.B

  $ws.Cells.Item($r.Row, 3).NumberFormat = "@"
  $ws.Cells.Item($r.Row, 3).Value = $r.C

  $ws.Cells.Item($r.Row, 4).Value = $r.D

  $ws.Cells.Item($r.Row, 5).NumberFormat = "@"
  $ws.Cells.Item($r.Row, 5).Value = $r.E

  $ws.Cells.Item($r.Row, 6).Value = $r.F

  $ws.Cells.Item($r.Row, 7).NumberFormat = "@"
  $ws.Cells.Item($r.Row, 7).Value = $r.G

  $ws.Cells.Item($r.Row, 8).NumberFormat = "@"
  $ws.Cells.Item($r.Row, 8).Value = $r.H
}

# --- 2. Restore the table's real formatting -------------------------------
# Row 98 already carries the exact style pattern used by every detail row:
#   A=s13 (left), B=(none), C=s13 (left), D=s14 (center), E=s14 (center),
#   F=s15 (center, red font), G=(none), H=(none)
# Paste its *formats only* onto each new row so every cell lands on the
# same style index without minting brand-new font/xf entries, then fix up
# the handful of F cells that should use the non-highlighted "s14" style
# (days-listed below the red-flag threshold) instead.
foreach ($r in $rows) {
  $ws.Range("A98:H98").Copy() | Out-Null
  $ws.Range("A" + $r.Row + ":H" + $r.Row).PasteSpecial(-4122) | Out-Null

  if ($r.FStyle -eq "s14") {
    $ws.Range("F92").Copy() | Out-Null
    $ws.Range("F" + $r.Row).PasteSpecial(-4122) | Out-Null
  }
}

$excel.CutCopyMode = 0
